$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values (PROPOSTA / RPS)
$ws.Range("B2").Value = 120339
$ws.Range("C2").Value = 12345

# Fill in row 3 (new data row)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 120340
$ws.Range("C3").Value = 12345

# Update the active selection to D3 as in the diff
$ws.Range("D3").Select()
